# Adapt column header formatting to respective input file names:
#   "<name>_old"  -> "<name>_FV2404"
#   "<name>_new"  -> "<name>_FV2410"
# and wrap the data range in an Excel Table (with autofilter), plus
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) cells ------------------------------
# Columns A-J carry the "_old" suffixed headers, column K is "diff"
# (left untouched), columns L-U carry the "_new" suffixed headers. Derive
# the replacement from whatever is already in the cell so this keeps
# working even if the exact header wording changes.
$oldSuffix = "_old"
$newSuffix = "_new"
$oldReplacement = "_FV2404"
$newReplacement = "_FV2410"

$lastCol = 21
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value()
    if ($header -like "*$oldSuffix") {
        $cell.Value = $header.Substring(0, $header.Length - $oldSuffix.Length) + $oldReplacement
    }
    elseif ($header -like "*$newSuffix") {
        $cell.Value = $header.Substring(0, $header.Length - $newSuffix.Length) + $newReplacement
    }
}

# --- 2. Freeze the header row --------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into a proper Excel Table --------------------
$dataRange = $ws.Range("A1:U59")
$table = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$table.Name = "Table1"

Write-Output "done"
